$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2023-11-19 Sunday" "2023-11-20 Monday"

Replace-Text "85÷5=" "71÷5="
Replace-Text "97÷3=" "83÷8="
Replace-Text "10÷9=" "93÷6="
Replace-Text "92÷8=" "93÷9="
Replace-Text "37÷8=" "95÷2="
Replace-Text "61÷8=" "57÷6="
Replace-Text "79÷2=" "62÷5="
Replace-Text "32÷2=" "58÷4="
Replace-Text "82÷2=" "25÷9="
Replace-Text "61÷9=" "68÷7="
Replace-Text "24÷3=" "35÷6="
Replace-Text "11÷2=" "56÷9="
Replace-Text "94÷2=" "79÷6="
Replace-Text "33÷5=" "56÷9="
Replace-Text "64÷7=" "95÷4="
Replace-Text "72÷3=" "41÷4="
Replace-Text "63÷9=" "67÷7="
Replace-Text "99÷7=" "77÷4="
Replace-Text "87÷6=" "79÷8="
Replace-Text "46÷6=" "98÷6="
Replace-Text "38÷4=" "65÷4="
Replace-Text "68÷4=" "91÷2="
Replace-Text "25÷6=" "96÷5="
Replace-Text "33÷2=" "31÷9="
Replace-Text "91÷5=" "43÷9="
